# Automatic update of files.
# Refreshes the "Förändrad" timestamp column, re-syncs the scraped per-case
# rows (Beteckning/Datum/Markägare/Area) and their dependent hyperlink
# formulas to the latest source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) date serial bump to 46074 for every data row (2-33)
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# 2) Per-row corrections: Beteckning (A), Datum (B), Markägare (F), Area ha (G)

# Row 4
$ws.Cells.Item(4, 1).Value = "A 389-2023"
$ws.Cells.Item(4, 2).Value = 44929
$ws.Cells.Item(4, 7).Value = 2.5

# Row 6
$ws.Cells.Item(6, 1).Value = "A 1782-2024"
$ws.Cells.Item(6, 2).Value = 45307
$ws.Cells.Item(6, 7).Value = 2.7

# Row 8
$ws.Cells.Item(8, 1).Value = "A 4481-2024"
$ws.Cells.Item(8, 2).Value = 45327
$ws.Cells.Item(8, 7).Value = 1

# Row 9
$ws.Cells.Item(9, 1).Value = "A 18327-2025"
$ws.Cells.Item(9, 2).Value = 45762
$ws.Cells.Item(9, 7).Value = 0.6

# Row 10
$ws.Cells.Item(10, 1).Value = "A 18328-2025"
$ws.Cells.Item(10, 2).Value = 45762
$ws.Cells.Item(10, 7).Value = 1.8

# Row 11
$ws.Cells.Item(11, 1).Value = "A 28260-2023"
$ws.Cells.Item(11, 2).Value = 45099
$ws.Cells.Item(11, 7).Value = 5

# Row 12 (also gains a Markägare entry)
$ws.Cells.Item(12, 1).Value = "A 10710-2025"
$ws.Cells.Item(12, 2).Value = 45722
$ws.Cells.Item(12, 6).Value = "Kommuner"
$ws.Cells.Item(12, 7).Value = 1.8

# Row 13
$ws.Cells.Item(13, 1).Value = "A 4822-2023"
$ws.Cells.Item(13, 2).Value = 44957
$ws.Cells.Item(13, 7).Value = 2.2

# Row 14 (loses its Markägare entry)
$ws.Cells.Item(14, 1).Value = "A 34400-2025"
$ws.Cells.Item(14, 2).Value = 45846.61351851852
$ws.Cells.Item(14, 6).Value = $null
$ws.Cells.Item(14, 7).Value = 1.3

# Row 15
$ws.Cells.Item(15, 1).Value = "A 34401-2025"
$ws.Cells.Item(15, 2).Value = 45846.6140162037
$ws.Cells.Item(15, 7).Value = 2.8

# Row 16
$ws.Cells.Item(16, 1).Value = "A 1531-2022"
$ws.Cells.Item(16, 2).Value = 44573
$ws.Cells.Item(16, 7).Value = 1.6

# Row 17
$ws.Cells.Item(17, 1).Value = "A 4256-2025"
$ws.Cells.Item(17, 2).Value = 45685
$ws.Cells.Item(17, 7).Value = 2

# Row 18
$ws.Cells.Item(18, 1).Value = "A 24-2023"
$ws.Cells.Item(18, 2).Value = 44928
$ws.Cells.Item(18, 7).Value = 0.5

# Row 19
$ws.Cells.Item(19, 1).Value = "A 11517-2024"
$ws.Cells.Item(19, 2).Value = 45372
$ws.Cells.Item(19, 7).Value = 0.7

# Row 20
$ws.Cells.Item(20, 1).Value = "A 4486-2024"
$ws.Cells.Item(20, 2).Value = 45327
$ws.Cells.Item(20, 7).Value = 0.6

# Row 21
$ws.Cells.Item(21, 1).Value = "A 21572-2023"
$ws.Cells.Item(21, 2).Value = 45063
$ws.Cells.Item(21, 7).Value = 1.7

# Row 22
$ws.Cells.Item(22, 1).Value = "A 635-2023"
$ws.Cells.Item(22, 2).Value = 44930
$ws.Cells.Item(22, 7).Value = 0.5

# Row 23
$ws.Cells.Item(23, 1).Value = "A 32610-2024"
$ws.Cells.Item(23, 2).Value = 45513
$ws.Cells.Item(23, 7).Value = 0.5

# Row 25
$ws.Cells.Item(25, 1).Value = "A 7731-2026"
$ws.Cells.Item(25, 2).Value = 46062.52008101852
$ws.Cells.Item(25, 7).Value = 5.9

# Row 26
$ws.Cells.Item(26, 1).Value = "A 18434-2023"
$ws.Cells.Item(26, 2).Value = 45042
$ws.Cells.Item(26, 7).Value = 0.7

# Row 27
$ws.Cells.Item(27, 1).Value = "A 4487-2024"
$ws.Cells.Item(27, 2).Value = 45327
$ws.Cells.Item(27, 7).Value = 1.9

# Row 28
$ws.Cells.Item(28, 1).Value = "A 5817-2025"
$ws.Cells.Item(28, 2).Value = 45694.74113425926
$ws.Cells.Item(28, 7).Value = 1.2

# Row 29
$ws.Cells.Item(29, 1).Value = "A 48974-2023"
$ws.Cells.Item(29, 2).Value = 45209
$ws.Cells.Item(29, 7).Value = 4.5

# Row 30
$ws.Cells.Item(30, 1).Value = "A 48181-2024"
$ws.Cells.Item(30, 2).Value = 45589

# Row 32
$ws.Cells.Item(32, 1).Value = "A 18332-2025"
$ws.Cells.Item(32, 2).Value = 45762
$ws.Cells.Item(32, 7).Value = 2.5

# Row 33
$ws.Cells.Item(33, 1).Value = "A 4493-2024"
$ws.Cells.Item(33, 2).Value = 45327
$ws.Cells.Item(33, 7).Value = 1.8

# 3) Rows 4 and 6 carry dependent HYPERLINK() formulas (artfynd/kartor/
#    klagomål/klagomålsmail/tillsyn/tillsynsmail) keyed off their Beteckning
#    (column A); regenerate them to match the swapped case numbers.
function Set-CaseLinks($row, $beteckning) {
    $ws.Cells.Item($row, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/artfynd/' + $beteckning + ' artfynd.xlsx", "' + $beteckning + '")'
    $ws.Cells.Item($row, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/kartor/' + $beteckning + ' karta.png", "' + $beteckning + '")'
    $ws.Cells.Item($row, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/klagomål/' + $beteckning + ' FSC-klagomål.docx", "' + $beteckning + '")'
    $ws.Cells.Item($row, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/klagomålsmail/' + $beteckning + ' FSC-klagomål mail.docx", "' + $beteckning + '")'
    $ws.Cells.Item($row, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/tillsyn/' + $beteckning + ' tillsynsbegäran.docx", "' + $beteckning + '")'
    $ws.Cells.Item($row, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1278/tillsynsmail/' + $beteckning + ' tillsynsbegäran mail.docx", "' + $beteckning + '")'
}

Set-CaseLinks 4 "A 389-2023"
Set-CaseLinks 6 "A 1782-2024"
